$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Block 1: rows 103-147 rotate up by one (CGB/Cuiaba moves to end)
$ws.Range('A103').Value2 = 'CWB'
$ws.Range('B103').Value2 = 'Curitiba, Brazil'
$ws.Range('C103').Value2 = 'South America'
$ws.Range('D103').Value2 = 'Curitiba'
$ws.Range('E103').Value2 = 'Brazil'
$ws.Range('F103').Value2 = 'BR'
$ws.Range('G103').Value2 = -25.5284996033
$ws.Range('H103').Value2 = -49.1758003235
$ws.Range('A104').Value2 = 'FLN'
$ws.Range('B104').Value2 = 'Florianopolis, Brazil'
$ws.Range('C104').Value2 = 'South America'
$ws.Range('D104').Value2 = 'Florianopolis'
$ws.Range('E104').Value2 = 'Brazil'
$ws.Range('F104').Value2 = 'BR'
$ws.Range('G104').Value2 = -27.6702785492
$ws.Range('H104').Value2 = -48.5525016785
$ws.Range('A105').Value2 = 'FOR'
$ws.Range('B105').Value2 = 'Fortaleza, Brazil'
$ws.Range('C105').Value2 = 'South America'
$ws.Range('D105').Value2 = 'Fortaleza'
$ws.Range('E105').Value2 = 'Brazil'
$ws.Range('F105').Value2 = 'BR'
$ws.Range('G105').Value2 = -3.7762799263
$ws.Range('H105').Value2 = -38.5326004028
$ws.Range('A106').Value2 = 'GEO'
$ws.Range('B106').Value2 = 'Georgetown, Guyana'
$ws.Range('C106').Value2 = 'South America'
$ws.Range('D106').Value2 = 'Georgetown'
$ws.Range('E106').Value2 = 'Guyana'
$ws.Range('F106').Value2 = 'GY'
$ws.Range('G106').Value2 = 6.825648
$ws.Range('H106').Value2 = -58.163756
$ws.Range('A107').Value2 = 'GYN'
$ws.Range('B107').Value2 = 'Goiania, Brazil'
$ws.Range('C107').Value2 = 'South America'
$ws.Range('D107').Value2 = 'Goiania'
$ws.Range('E107').Value2 = 'Brazil'
$ws.Range('F107').Value2 = 'BR'
$ws.Range('G107').Value2 = -16.69727
$ws.Range('H107').Value2 = -49.26851
$ws.Range('A108').Value2 = 'GUA'
$ws.Range('B108').Value2 = 'Guatemala City, Guatemala'
$ws.Range('C108').Value2 = 'North America'
$ws.Range('D108').Value2 = 'Guatemala City'
$ws.Range('E108').Value2 = 'Guatemala'
$ws.Range('F108').Value2 = 'GT'
$ws.Range('G108').Value2 = 14.5832996368
$ws.Range('H108').Value2 = -90.5274963379
$ws.Range('A109').Value2 = 'GYE'
$ws.Range('B109').Value2 = 'Guayaquil, Ecuador'
$ws.Range('C109').Value2 = 'South America'
$ws.Range('D109').Value2 = 'Guayaquil'
$ws.Range('E109').Value2 = 'Ecuador'
$ws.Range('F109').Value2 = 'EC'
$ws.Range('G109').Value2 = -2.1894
$ws.Range('H109').Value2 = -79.8891
$ws.Range('A110').Value2 = 'ITJ'
$ws.Range('B110').Value2 = 'Itajai, Brazil'
$ws.Range('C110').Value2 = 'South America'
$ws.Range('D110').Value2 = 'Itajai'
$ws.Range('E110').Value2 = 'Brazil'
$ws.Range('F110').Value2 = 'BR'
$ws.Range('G110').Value2 = -27.6116676331
$ws.Range('H110').Value2 = -48.6727790833
$ws.Range('A111').Value2 = 'JOI'
$ws.Range('B111').Value2 = 'Joinville, Brazil'
$ws.Range('C111').Value2 = 'South America'
$ws.Range('D111').Value2 = 'Joinville'
$ws.Range('E111').Value2 = 'Brazil'
$ws.Range('F111').Value2 = 'BR'
$ws.Range('G111').Value2 = -26.304408
$ws.Range('H111').Value2 = -48.846383
$ws.Range('A112').Value2 = 'JDO'
$ws.Range('B112').Value2 = 'Juazeiro do Norte, Brazil'
$ws.Range('C112').Value2 = 'South America'
$ws.Range('D112').Value2 = 'Juazeiro do Norte'
$ws.Range('E112').Value2 = 'Brazil'
$ws.Range('F112').Value2 = 'BR'
$ws.Range('G112').Value2 = -7.2242
$ws.Range('H112').Value2 = -39.313
$ws.Range('A113').Value2 = 'LIM'
$ws.Range('B113').Value2 = 'Lima, Peru'
$ws.Range('C113').Value2 = 'South America'
$ws.Range('D113').Value2 = 'Lima'
$ws.Range('E113').Value2 = 'Peru'
$ws.Range('F113').Value2 = 'PE'
$ws.Range('G113').Value2 = -12.021900177
$ws.Range('H113').Value2 = -77.1143035889
$ws.Range('A114').Value2 = 'MAO'
$ws.Range('B114').Value2 = 'Manaus, Brazil'
$ws.Range('C114').Value2 = 'South America'
$ws.Range('D114').Value2 = 'Manaus'
$ws.Range('E114').Value2 = 'Brazil'
$ws.Range('F114').Value2 = 'BR'
$ws.Range('G114').Value2 = -3.11286
$ws.Range('H114').Value2 = -60.01949
$ws.Range('A115').Value2 = 'MDE'
$ws.Range('B115').Value2 = 'Medellín, Colombia'
$ws.Range('C115').Value2 = 'South America'
$ws.Range('D115').Value2 = 'Medellín'
$ws.Range('E115').Value2 = 'Colombia'
$ws.Range('F115').Value2 = 'CO'
$ws.Range('G115').Value2 = 6.16454
$ws.Range('H115').Value2 = -75.42310000000001
$ws.Range('A116').Value2 = 'NQN'
$ws.Range('B116').Value2 = 'Neuquen, Argentina'
$ws.Range('C116').Value2 = 'South America'
$ws.Range('D116').Value2 = 'Neuquen'
$ws.Range('E116').Value2 = 'Argentina'
$ws.Range('F116').Value2 = 'AR'
$ws.Range('G116').Value2 = -38.9490013123
$ws.Range('H116').Value2 = -68.1557006836
$ws.Range('A117').Value2 = 'PTY'
$ws.Range('B117').Value2 = 'Panama City, Panama'
$ws.Range('C117').Value2 = 'South America'
$ws.Range('D117').Value2 = 'Panama City'
$ws.Range('E117').Value2 = 'Panama'
$ws.Range('F117').Value2 = 'PA'
$ws.Range('G117').Value2 = 9.0713596344
$ws.Range('H117').Value2 = -79.3834991455
$ws.Range('A118').Value2 = 'PBM'
$ws.Range('B118').Value2 = 'Paramaribo, Suriname'
$ws.Range('C118').Value2 = 'South America'
$ws.Range('D118').Value2 = 'Paramaribo'
$ws.Range('E118').Value2 = 'Suriname'
$ws.Range('F118').Value2 = 'SR'
$ws.Range('G118').Value2 = 5.452831
$ws.Range('H118').Value2 = -55.187783
$ws.Range('A119').Value2 = 'POA'
$ws.Range('B119').Value2 = 'Porto Alegre, Brazil'
$ws.Range('C119').Value2 = 'South America'
$ws.Range('D119').Value2 = 'Porto Alegre'
$ws.Range('E119').Value2 = 'Brazil'
$ws.Range('F119').Value2 = 'BR'
$ws.Range('G119').Value2 = -29.9944000244
$ws.Range('H119').Value2 = -51.1713981628
$ws.Range('A120').Value2 = 'UIO'
$ws.Range('B120').Value2 = 'Quito, Ecuador'
$ws.Range('C120').Value2 = 'South America'
$ws.Range('D120').Value2 = 'Quito'
$ws.Range('E120').Value2 = 'Ecuador'
$ws.Range('F120').Value2 = 'EC'
$ws.Range('G120').Value2 = -0.1291666667
$ws.Range('H120').Value2 = -78.3575
$ws.Range('A121').Value2 = 'REC'
$ws.Range('B121').Value2 = 'Recife, Brazil'
$ws.Range('C121').Value2 = 'South America'
$ws.Range('D121').Value2 = 'Recife'
$ws.Range('E121').Value2 = 'Brazil'
$ws.Range('F121').Value2 = 'BR'
$ws.Range('G121').Value2 = -8.126489639300001
$ws.Range('H121').Value2 = -34.9235992432
$ws.Range('A122').Value2 = 'RAO'
$ws.Range('B122').Value2 = 'Ribeirao Preto, Brazil'
$ws.Range('C122').Value2 = 'South America'
$ws.Range('D122').Value2 = 'Ribeirao Preto'
$ws.Range('E122').Value2 = 'Brazil'
$ws.Range('F122').Value2 = 'BR'
$ws.Range('G122').Value2 = -21.1363887787
$ws.Range('H122').Value2 = -47.7766685486
$ws.Range('A123').Value2 = 'GIG'
$ws.Range('B123').Value2 = 'Rio de Janeiro, Brazil'
$ws.Range('C123').Value2 = 'South America'
$ws.Range('D123').Value2 = 'Rio de Janeiro'
$ws.Range('E123').Value2 = 'Brazil'
$ws.Range('F123').Value2 = 'BR'
$ws.Range('G123').Value2 = -22.8099994659
$ws.Range('H123').Value2 = -43.2505569458
$ws.Range('A124').Value2 = 'SJO'
$ws.Range('B124').Value2 = 'San José, Costa Rica'
$ws.Range('C124').Value2 = 'South America'
$ws.Range('D124').Value2 = 'San José'
$ws.Range('E124').Value2 = 'Costa Rica'
$ws.Range('F124').Value2 = 'CR'
$ws.Range('G124').Value2 = 9.9938602448
$ws.Range('H124').Value2 = -84.2088012695
$ws.Range('A125').Value2 = 'SCL'
$ws.Range('B125').Value2 = 'Santiago, Chile'
$ws.Range('C125').Value2 = 'South America'
$ws.Range('D125').Value2 = 'Santiago'
$ws.Range('E125').Value2 = 'Chile'
$ws.Range('F125').Value2 = 'CL'
$ws.Range('G125').Value2 = -33.3930015564
$ws.Range('H125').Value2 = -70.7857971191
$ws.Range('A126').Value2 = 'SDQ'
$ws.Range('B126').Value2 = 'Santo Domingo, Dominican Republic'
$ws.Range('C126').Value2 = 'North America'
$ws.Range('D126').Value2 = 'Santo Domingo'
$ws.Range('E126').Value2 = 'Dominican Republic'
$ws.Range('F126').Value2 = 'DO'
$ws.Range('G126').Value2 = 18.4297008514
$ws.Range('H126').Value2 = -69.6688995361
$ws.Range('A127').Value2 = 'SJP'
$ws.Range('B127').Value2 = 'São José do Rio Preto, Brazil'
$ws.Range('C127').Value2 = 'South America'
$ws.Range('D127').Value2 = 'São José do Rio Preto'
$ws.Range('E127').Value2 = 'Brazil'
$ws.Range('F127').Value2 = 'BR'
$ws.Range('G127').Value2 = -20.807157
$ws.Range('H127').Value2 = -49.378994
$ws.Range('A128').Value2 = 'SJK'
$ws.Range('B128').Value2 = 'São José dos Campos, Brazil'
$ws.Range('C128').Value2 = 'South America'
$ws.Range('D128').Value2 = 'São José dos Campos'
$ws.Range('E128').Value2 = 'Brazil'
$ws.Range('F128').Value2 = 'BR'
$ws.Range('G128').Value2 = -23.1791
$ws.Range('H128').Value2 = -45.8872
$ws.Range('A129').Value2 = 'GRU'
$ws.Range('B129').Value2 = 'São Paulo, Brazil'
$ws.Range('C129').Value2 = 'South America'
$ws.Range('D129').Value2 = 'São Paulo'
$ws.Range('E129').Value2 = 'Brazil'
$ws.Range('F129').Value2 = 'BR'
$ws.Range('G129').Value2 = -23.4355564117
$ws.Range('H129').Value2 = -46.4730567932
$ws.Range('A130').Value2 = 'SOD'
$ws.Range('B130').Value2 = 'Sorocaba, Brazil'
$ws.Range('C130').Value2 = 'South America'
$ws.Range('D130').Value2 = 'Sorocaba'
$ws.Range('E130').Value2 = 'Brazil'
$ws.Range('F130').Value2 = 'BR'
$ws.Range('G130').Value2 = -23.54389
$ws.Range('H130').Value2 = -46.63445
$ws.Range('A131').Value2 = 'TGU'
$ws.Range('B131').Value2 = 'Tegucigalpa, Honduras'
$ws.Range('C131').Value2 = 'South America'
$ws.Range('D131').Value2 = 'Tegucigalpa'
$ws.Range('E131').Value2 = 'Honduras'
$ws.Range('F131').Value2 = 'HN'
$ws.Range('G131').Value2 = 14.0608
$ws.Range('H131').Value2 = -87.21720000000001
$ws.Range('A132').Value2 = 'NVT'
$ws.Range('B132').Value2 = 'Timbo, Brazil'
$ws.Range('C132').Value2 = 'South America'
$ws.Range('D132').Value2 = 'Timbo'
$ws.Range('E132').Value2 = 'Brazil'
$ws.Range('F132').Value2 = 'BR'
$ws.Range('G132').Value2 = -26.8251
$ws.Range('H132').Value2 = -49.2695
$ws.Range('A133').Value2 = 'UDI'
$ws.Range('B133').Value2 = 'Uberlandia, Brazil'
$ws.Range('C133').Value2 = 'South America'
$ws.Range('D133').Value2 = 'Uberlandia'
$ws.Range('E133').Value2 = 'Brazil'
$ws.Range('F133').Value2 = 'BR'
$ws.Range('G133').Value2 = -18.8836116791
$ws.Range('H133').Value2 = -48.225276947
$ws.Range('A134').Value2 = 'VIX'
$ws.Range('B134').Value2 = 'Vitoria, Brazil'
$ws.Range('C134').Value2 = 'South America'
$ws.Range('D134').Value2 = 'Vitoria'
$ws.Range('E134').Value2 = 'Brazil'
$ws.Range('F134').Value2 = 'BR'
$ws.Range('G134').Value2 = -20.64871
$ws.Range('H134').Value2 = -41.90857
$ws.Range('A135').Value2 = 'CAW'
$ws.Range('B135').Value2 = 'Campos dos Goytacazes, Brazil'
$ws.Range('C135').Value2 = 'South America'
$ws.Range('D135').Value2 = 'Campos dos Goytacazes'
$ws.Range('E135').Value2 = 'Brazil'
$ws.Range('F135').Value2 = 'BR'
$ws.Range('G135').Value2 = -21.698299408
$ws.Range('H135').Value2 = -41.301700592
$ws.Range('A136').Value2 = 'XAP'
$ws.Range('B136').Value2 = 'Chapeco, Brazil'
$ws.Range('C136').Value2 = 'South America'
$ws.Range('D136').Value2 = 'Chapeco'
$ws.Range('E136').Value2 = 'Brazil'
$ws.Range('F136').Value2 = 'BR'
$ws.Range('G136').Value2 = -27.1341991425
$ws.Range('H136').Value2 = -52.6566009521
$ws.Range('A137').Value2 = 'BGI'
$ws.Range('B137').Value2 = 'Bridgetown, Barbados'
$ws.Range('C137').Value2 = 'North America'
$ws.Range('D137').Value2 = 'Bridgetown'
$ws.Range('E137').Value2 = 'Barbados'
$ws.Range('F137').Value2 = 'BB'
$ws.Range('G137').Value2 = 13.103562
$ws.Range('H137').Value2 = -59.603226
$ws.Range('A138').Value2 = 'GND'
$ws.Range('B138').Value2 = 'St. George''s, Grenada'
$ws.Range('C138').Value2 = 'South America'
$ws.Range('D138').Value2 = 'St. George''s'
$ws.Range('E138').Value2 = 'Grenada'
$ws.Range('F138').Value2 = 'GD'
$ws.Range('G138').Value2 = 12.007116
$ws.Range('H138').Value2 = -61.7882288
$ws.Range('A139').Value2 = 'STI'
$ws.Range('B139').Value2 = 'Santiago de los Caballeros, Dominican Republic'
$ws.Range('C139').Value2 = 'North America'
$ws.Range('D139').Value2 = 'Santiago de los Caballeros'
$ws.Range('E139').Value2 = 'Dominican Republic'
$ws.Range('F139').Value2 = 'DO'
$ws.Range('G139').Value2 = 19.4060993195
$ws.Range('H139').Value2 = -70.60469818120001
$ws.Range('A140').Value2 = 'LPB'
$ws.Range('B140').Value2 = 'La Paz, Bolivia'
$ws.Range('C140').Value2 = 'South America'
$ws.Range('D140').Value2 = 'La Paz'
$ws.Range('E140').Value2 = 'Bolivia'
$ws.Range('F140').Value2 = 'BO'
$ws.Range('G140').Value2 = -16.4897
$ws.Range('H140').Value2 = -68.1193
$ws.Range('A141').Value2 = 'SJU'
$ws.Range('B141').Value2 = 'San Juan, Puerto Rico'
$ws.Range('C141').Value2 = 'North America'
$ws.Range('D141').Value2 = 'San Juan'
$ws.Range('E141').Value2 = 'Puerto Rico'
$ws.Range('F141').Value2 = 'PR'
$ws.Range('G141').Value2 = 18.411391
$ws.Range('H141').Value2 = -66.10279300000001
$ws.Range('A142').Value2 = 'BAQ'
$ws.Range('B142').Value2 = 'Barranquilla, Colombia'
$ws.Range('C142').Value2 = 'South America'
$ws.Range('D142').Value2 = 'Barranquilla'
$ws.Range('E142').Value2 = 'Colombia'
$ws.Range('F142').Value2 = 'CO'
$ws.Range('G142').Value2 = 10.8896
$ws.Range('H142').Value2 = -74.7808
$ws.Range('A143').Value2 = 'PMW'
$ws.Range('B143').Value2 = 'Palmas, Brazil'
$ws.Range('C143').Value2 = 'South America'
$ws.Range('D143').Value2 = 'Palmas'
$ws.Range('E143').Value2 = 'Brazil'
$ws.Range('F143').Value2 = 'BR'
$ws.Range('G143').Value2 = -10.2915000916
$ws.Range('H143').Value2 = -48.3569984436
$ws.Range('A144').Value2 = 'ARU'
$ws.Range('B144').Value2 = 'Aracatuba, Brazil'
$ws.Range('C144').Value2 = 'South America'
$ws.Range('D144').Value2 = 'Aracatuba'
$ws.Range('E144').Value2 = 'Brazil'
$ws.Range('F144').Value2 = 'BR'
$ws.Range('G144').Value2 = -21.1413002014
$ws.Range('H144').Value2 = -50.4247016907
$ws.Range('A145').Value2 = 'POS'
$ws.Range('B145').Value2 = 'Port of Spain, Trinidad and Tobago'
$ws.Range('C145').Value2 = 'South America'
$ws.Range('D145').Value2 = 'Port of Spain'
$ws.Range('E145').Value2 = 'Trinidad and Tobago'
$ws.Range('F145').Value2 = 'TT'
$ws.Range('G145').Value2 = 10.5953998566
$ws.Range('H145').Value2 = -61.3372001648
$ws.Range('A146').Value2 = 'SSA'
$ws.Range('B146').Value2 = 'Salvador, Brazil'
$ws.Range('C146').Value2 = 'South America'
$ws.Range('D146').Value2 = 'Salvador'
$ws.Range('E146').Value2 = 'Brazil'
$ws.Range('F146').Value2 = 'BR'
$ws.Range('G146').Value2 = -12.9086112976
$ws.Range('H146').Value2 = -38.3224983215
$ws.Range('A147').Value2 = 'CGB'
$ws.Range('B147').Value2 = 'Cuiaba, Brazil'
$ws.Range('C147').Value2 = 'South America'
$ws.Range('D147').Value2 = 'Cuiaba'
$ws.Range('E147').Value2 = 'Brazil'
$ws.Range('F147').Value2 = 'BR'
$ws.Range('G147').Value2 = -15.59611
$ws.Range('H147').Value2 = -56.09667

# Block 2: rows 172-179 rotate up by one (RUN/Saint-Denis moves to end)
$ws.Range('A172').Value2 = 'TUN'
$ws.Range('B172').Value2 = 'Tunis, Tunisia'
$ws.Range('C172').Value2 = 'Africa'
$ws.Range('D172').Value2 = 'Tunis'
$ws.Range('E172').Value2 = 'Tunisia'
$ws.Range('F172').Value2 = 'TN'
$ws.Range('G172').Value2 = 36.8510017395
$ws.Range('H172').Value2 = 10.2271995544
$ws.Range('A173').Value2 = 'FIH'
$ws.Range('B173').Value2 = 'Kinshasa, DR Congo'
$ws.Range('C173').Value2 = 'Africa'
$ws.Range('D173').Value2 = 'Kinshasa'
$ws.Range('E173').Value2 = 'DR Congo'
$ws.Range('F173').Value2 = 'CD'
$ws.Range('G173').Value2 = -4.3857498169
$ws.Range('H173').Value2 = 15.4446001053
$ws.Range('A174').Value2 = 'CAI'
$ws.Range('B174').Value2 = 'Cairo, Egypt'
$ws.Range('C174').Value2 = 'Africa'
$ws.Range('D174').Value2 = 'Cairo'
$ws.Range('E174').Value2 = 'Egypt'
$ws.Range('F174').Value2 = 'EG'
$ws.Range('G174').Value2 = 30.1219005585
$ws.Range('H174').Value2 = 31.4055995941
$ws.Range('A175').Value2 = 'WDH'
$ws.Range('B175').Value2 = 'Windhoek, Namibia'
$ws.Range('C175').Value2 = 'Africa'
$ws.Range('D175').Value2 = 'Windhoek'
$ws.Range('E175').Value2 = 'Namibia'
$ws.Range('F175').Value2 = 'NA'
$ws.Range('G175').Value2 = -22.565587
$ws.Range('H175').Value2 = 17.085334
$ws.Range('A176').Value2 = 'ASK'
$ws.Range('B176').Value2 = 'Yamoussoukro, Ivory Coast'
$ws.Range('C176').Value2 = 'Africa'
$ws.Range('D176').Value2 = 'Yamoussoukro'
$ws.Range('E176').Value2 = 'Ivory Coast'
$ws.Range('F176').Value2 = 'CI'
$ws.Range('G176').Value2 = 6.842178
$ws.Range('H176').Value2 = -5.259932
$ws.Range('A177').Value2 = 'ABJ'
$ws.Range('B177').Value2 = 'Abidjan, Ivory Coast'
$ws.Range('C177').Value2 = 'Africa'
$ws.Range('D177').Value2 = 'Abidjan'
$ws.Range('E177').Value2 = 'Ivory Coast'
$ws.Range('F177').Value2 = 'CI'
$ws.Range('G177').Value2 = 5.292598
$ws.Range('H177').Value2 = -3.999133
$ws.Range('A178').Value2 = 'EBB'
$ws.Range('B178').Value2 = 'Kampala, Uganda'
$ws.Range('C178').Value2 = 'Africa'
$ws.Range('D178').Value2 = 'Kampala'
$ws.Range('E178').Value2 = 'Uganda'
$ws.Range('F178').Value2 = 'UG'
$ws.Range('G178').Value2 = 0.3152
$ws.Range('H178').Value2 = 32.5816
$ws.Range('A179').Value2 = 'RUN'
$ws.Range('B179').Value2 = 'Saint-Denis, Réunion'
$ws.Range('C179').Value2 = 'Africa'
$ws.Range('D179').Value2 = 'Saint-Denis'
$ws.Range('E179').Value2 = 'Réunion'
$ws.Range('F179').Value2 = 'RE'
$ws.Range('G179').Value2 = -20.8871002197
$ws.Range('H179').Value2 = 55.5102996826

# Row 310: update name/country text for Sioux Falls
$ws.Range('B310').Value2 = 'Sioux Falls, SD, United States'
$ws.Range('E310').Value2 = 'United States'

